# pinout.xlsx update: "Remove freertos, new pinout"
# Re-labels several pins, adds GPIO alt-function columns (H / I) to the
# Senzory table, removes now-unused PTA4/FTM0_CH2/Pull Up-Down strings,
# clears the leftover column-A marker cells, and tidies up the "Hallovky"
# and last "OpenSDA" tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------
# "Hallovky" table (rows 10-11): swap PTD3 -> PTD2 for "Right hall" and
# add an IRQ indicator column (H) mirroring column G for both rows.
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "PTD2"

Copy-Format "G10" "H10"
$ws.Range("H10").Value = "IRQ"

Copy-Format "G11" "H11"
$ws.Range("H11").Value = "IRQ"

# ---------------------------------------------------------------------
# "Senzory" table (rows 15-21): pin/alt-func relabeling.
# ---------------------------------------------------------------------

# Row 15 note moved from column H to column I (same formatting).
$ws.Range("I15").Value = $ws.Range("H15").Value
Copy-Format "H15" "I15"
$ws.Range("H15").Clear() | Out-Null

# Row 16: PTD1 -> PTD4, FTM0_CH1 -> FTM0_CH4
$ws.Range("B16").Value = "PTD4"
$ws.Range("D16").Value = "FTM0_CH4"

# Row 17: PTD2 -> PTD3, FTM0_CH2 -> FTM0_CH3
$ws.Range("B17").Value = "PTD3"
$ws.Range("D17").Value = "FTM0_CH3"

# Row 18: PTA4 -> PTA0, "Pull Up/Down" -> "GPIOA, 0"; add IRQ in column H.
$ws.Range("B18").Value = "PTA0"
$ws.Range("D18").Value = "GPIOA, 0"
Copy-Format "G18" "H18"
$ws.Range("H18").Value = "IRQ"

# Row 19: "Pull Up/Down" -> "GPIOA, 5"; add IRQ in column H.
$ws.Range("D19").Value = "GPIOA, 5"
Copy-Format "G19" "H19"
$ws.Range("H19").Value = "IRQ"

# Row 20: add "GPIOC, 4" device id and a GPIO column entry.
$ws.Range("D20").Value = "GPIOC, 4"
Copy-Format "C20" "H20"
$ws.Range("H20").Value = "GPIO"

# Row 21: add "GPIOC, 5" device id and a GPIO column entry (matches the
# source workbook's own H20 styling, reused for H21 as well).
$ws.Range("D21").Value = "GPIOC, 5"
Copy-Format "C20" "H21"
$ws.Range("H21").Value = "GPIO"

# ---------------------------------------------------------------------
# Remove the leftover column-A marker cells next to the Senzory table.
# ---------------------------------------------------------------------
$ws.Range("A15:A17").Clear() | Out-Null

# ---------------------------------------------------------------------
# "OpenSDA" table (rows 30-31): PTA0 -> PTA2 for the TX pin.
# ---------------------------------------------------------------------
$ws.Range("B31").Value = "PTA2"

# ---------------------------------------------------------------------
# New PWM column (H) for the Motory table (rows 4-6).
# ---------------------------------------------------------------------
Copy-Format "G4" "H4"
$ws.Range("H4").Value = "PWM"

Copy-Format "G5" "H5"
$ws.Range("H5").Value = "PWM"

Copy-Format "G6" "H6"
$ws.Range("H6").Value = "PWM"

# ---------------------------------------------------------------------
# Column widths to match the narrow helper columns A and H.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.140625
$ws.Columns.Item(8).ColumnWidth = 9.140625

# ---------------------------------------------------------------------
# Selection / scroll position.
# ---------------------------------------------------------------------
$ws.Range("D18").Select() | Out-Null
